$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 190
$ws.Range("I2").Value = 488
$ws.Range("J2").Value = 2134
$ws.Range("K2").Value = 7
$ws.Range("L2").Value = 563
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 340
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 6
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 46
$ws.Range("S2").Value = 236
$ws.Range("T2").Value = 386
$ws.Range("U2").Value = 31
$ws.Range("V2").Value = 3227
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 3218
$ws.Range("Z2").Value = 47
$ws.Range("AA2").Value = 19
